$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill newly added day columns (L-P for rows 2-28, J-N for rows 29-36)
# matching the extended survival-tracking days added in this data update.

$ws.Range("L2:P4").Value = 0
$ws.Range("L5:P22").Value = 1
$ws.Range("L23:P24").Value = 0
$ws.Range("L26:P27").Value = 0
$ws.Range("J29:N34").Value = 1
$ws.Range("J35:N36").Value = 0

# Row 25 and Row 28 have non-uniform values across the new columns
$ws.Range("L25").Value = 1
$ws.Range("M25:P25").Value = 0
$ws.Range("L28:O28").Value = 1
$ws.Range("P28").Value = 0

# Update the view state (scrolled/selected cell) to match the saved session
$ws.Range("Q32").Select()
$excel.ActiveWindow.ScrollColumn = 8
